$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.287.93'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '3.536.32'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''608.04'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '''143.80'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("D7").Value = '3.535.12'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("E10").Value = '  -4.28%  '
$ws.Range("D11").Value = '''8.06'
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").Value = '4.135.68'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("E14").Value = '  -4.62%  '
$ws.Range("D15").Value = '''30.24'
$ws.Range("E15").Value = '  -5.34%  '
$ws.Range("D16").Value = '3.532.71'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '66.371.22'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").Value = '''10.93'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '''14.93'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("D22").Value = '''425.70'
$ws.Range("E22").Value = '  -2.87%  '
$ws.Range("D23").Value = '''0.601'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").Value = '''78.66'
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").Value = '3.677.27'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").Value = '''8.10'
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("D29").Value = '''9.18'
$ws.Range("E29").Value = '  -6.20%  '
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -8.02%  '
$ws.Range("D33").Value = '''0.160'
$ws.Range("E33").Value = '  -4.41%  '
$ws.Range("D34").Value = '''25.27'
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("D35").Value = '3.526.19'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '''1.76'
$ws.Range("E37").Value = '  -3.12%  '
$ws.Range("D38").Value = '''5.64'
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '''172.14'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("E42").Value = '  -4.45%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = '''1.90'
$ws.Range("E45").Value = '  -8.03%  '
$ws.Range("D46").Value = '''45.56'
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = '''26.14'
$ws.Range("E47").Value = '  -6.36%  '
$ws.Range("D48").Value = '''1.22'
$ws.Range("E48").Value = '  -4.52%  '
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("E50").Value = '  -4.40%  '
$ws.Range("E51").Value = '  -4.72%  '
